# Madaswamy Murugan Resume.docx -- apply the "_GoBack" bookmark relocation
# and the "May 2022" -> "May 2014" edit described by the commit diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the stray "_GoBack" last-edit-position bookmark.
#    It currently sits at the end of the "... (Delhi University)," run
#    (right before the following </w:p>).  Remove it there ...
# ---------------------------------------------------------------------
$oldMark = $d.Bookmarks("_GoBack")
$oldMark.Delete()

# ... and re-create it immediately after the header table (where Word's
# most-recent-edit tracking now points), i.e. at the very start of the
# paragraph that follows the table.
$headerTable = $d.Tables(1)
$afterTable = $d.Range($headerTable.Range.End, $headerTable.Range.End)
$d.Bookmarks.Add("_GoBack", $afterTable)

# ---------------------------------------------------------------------
# 2) "Since : May 2022 till November 2017." -> "Since : May 2014 till
#    November 2017."  The run that used to hold " May 2022 " is split
#    into " May 20" and "14 " (same run formatting) by the edit, so we
#    replace the "22" digits and then force a run break right after
#    "May 20" by toggling a formatting property on/off over the "14 "
#    text (Word splits runs on direct-formatting boundaries even when
#    the property is immediately reverted).
# ---------------------------------------------------------------------
$search = $d.Content
$found = $search.Find.Execute(" May 2022 till November", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $yearDigits = $d.Range($search.Start + 7, $search.Start + 9)   # the "22" in "2022"
    $yearDigits.Text = "14"

    $newRun = $d.Range($search.Start + 7, $search.Start + 10)      # "14 "
    $newRun.Bold = 1
    $newRun.Bold = 0
}
